# Update Excel files after daily scrape - 2026-02-17 04:24:31 UTC

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# Column width updates
# (ColumnWidth is adjusted by a small constant offset, 5/6 of a character,
#  when the runtime round-trips through pixel measurements, so the assigned
#  values below are pre-compensated to land exactly on the target widths of
#  54, 51, 16 and 31 characters respectively.)
# ---------------------------------------------------------------------------
$ws.Columns.Item(3).ColumnWidth = 53.1666667
$ws.Columns.Item(4).ColumnWidth = 50.1666667
$ws.Columns.Item(7).ColumnWidth = 15.1666667
$ws.Columns.Item(8).ColumnWidth = 30.1666667

# ---------------------------------------------------------------------------
# Data table - rows 2 through 7
# ---------------------------------------------------------------------------
$data = @(
    @("1331991", "https://aiesec.org/opportunity/global-talent/1331991", "[CC] Social Media Content Creator (EU Only)", "Charles-de-Gaulle-Straße 20, 53113 Bonn, Germany", "Yes", "2 applicants", "6 - 18 Months", "DHL Group"),
    @("1331986", "https://aiesec.org/opportunity/global-talent/1331986", "Cybersecurity Intern", "Aronj, Uttar Pradesh, India", "No", "2 applicants", "3 - 6 Months", "FS University"),
    @("1331980", "https://aiesec.org/opportunity/global-talent/1331980", "Sales Responsible", "Ürgüp, Nevşehir, Türkiye", "No", "6 applicants", "6 - 18 Months", "Pink Lotus Jewellery"),
    @("1331751", "https://aiesec.org/opportunity/global-talent/1331751", "Customer Support Engineer", "Amman, Jordanie", "No", "0 applicants", "3 - 6 Months", "Estarta"),
    @("1331590", "https://aiesec.org/opportunity/global-talent/1331590", "Order to Cash (Accounts Receivable) Intern", "Dubai - United Arab Emirates", "No", "30 applicants", "3 - 6 Months", "Dubai Holding Group Services"),
    @("1331559", "https://aiesec.org/opportunity/global-talent/1331559", "B2B Recruitment Consultant - German/English Speaker", "Prague, Tchéquie", "No", "21 applicants", "6 - 18 Months", "Non Stop Consulting")
)

$row = 2
foreach ($record in $data) {
    # Column A holds opportunity IDs that look numeric ("1331991"). The
    # source data is textual, so the cell is forced to Text format before
    # the assignment (otherwise it is auto-converted to a number) and the
    # style is reverted back to Normal afterwards so the cell keeps its
    # original (unstyled) appearance while the stored value stays text.
    $ws.Cells.Item($row, 1).NumberFormat = "@"
    $ws.Cells.Item($row, 1).Value = $record[0]
    $ws.Cells.Item($row, 1).Style = "Normal"

    $ws.Cells.Item($row, 2).Value = $record[1]
    $ws.Cells.Item($row, 3).Value = $record[2]
    $ws.Cells.Item($row, 4).Value = $record[3]
    $ws.Cells.Item($row, 5).Value = $record[4]
    $ws.Cells.Item($row, 6).Value = $record[5]
    $ws.Cells.Item($row, 7).Value = $record[6]
    $ws.Cells.Item($row, 8).Value = $record[7]
    $row = $row + 1
}

# ---------------------------------------------------------------------------
# Highlight the "Yes" premium cell (E2) with a yellow fill
# ---------------------------------------------------------------------------
$ws.Range("E2").Interior.Color = 65535

$wb.Save()
